# Actualización automática del mapa (2025-08-06 07:05:49)
# Append 3 new incident rows (81-83) to the "AYKO" sheet, mirroring the
# existing table layout: columns A-P, text columns stored as text and
# the Attachments (I) / coordinate (M,N) columns stored as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that must be stored as literal text even though some of the
# values look numeric (Caso, F. De Reclamo, Direccion, Comuna, OT, ...).
$textCols = 1,2,3,4,5,6,7,8,10,11,12,15,16

$rows = @(
    @{
        1  = "6193"
        2  = "8/5/2025"
        3  = "POLA 591"
        4  = "9"
        5  = "808720861"
        6  = "AYKO"
        7  = "Pendiente"
        8  = "Picada"
        9  = 0
        10 = "Cambio"
        11 = "Sin equipos"
        12 = "Pasante"
        13 = -58.507385
        14 = -34.644479
        15 = "Devoto"
        16 = "Capital Norte"
    },
    @{
        1  = "6277"
        2  = "8/5/2025"
        3  = "SENILLOSA 323"
        4  = "6"
        5  = "808720859"
        6  = "AYKO"
        7  = "Pendiente"
        8  = "Picada"
        9  = 1
        10 = "Cambio"
        11 = "Sin equipos"
        12 = "Pasante"
        13 = -58.429726
        14 = -34.619969
        15 = "Boedo"
        16 = "Capital Sur"
    },
    @{
        1  = "6274"
        2  = "8/5/2025"
        3  = "ARANGUREN, JUAN F., Dr. 2964"
        4  = "12"
        5  = "808720858"
        6  = "AYKO"
        7  = "Pendiente"
        8  = "chocada"
        9  = 1
        10 = "Cambio"
        11 = "Sin equipos"
        12 = "Pasante"
        13 = -58.473089
        14 = -34.625478
        15 = "Devoto"
        16 = "Capital Norte"
    }
)

$startRow = 81
$r = $startRow
foreach ($row in $rows) {
    foreach ($col in 1..16) {
        $cell = $ws.Cells.Item($r, $col)
        $value = $row[$col]
        if ($textCols -contains $col) {
            # Force text storage (matches the source data's inline-string
            # cells) instead of Excel's automatic number/date detection.
            $cell.NumberFormat = "@"
            $cell.Value = [string]$value
        } else {
            $cell.Value = $value
        }
    }
    $r++
}
